$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "function" header in D1, matching the bold/bordered/centered style
# already applied to B1/C1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "function"

# Helper cell used to coerce numeric-looking text into real text (t="s")
# values without inheriting any special number format/style: we stage the
# text in a scratch cell with NumberFormat "@" then paste-values it into
# the destination, which keeps the destination's existing (default) style.
$scratch = $ws.Range("Z1")

function Set-TextValue($cell, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

$centers = @("588.2324014377622", "590.579897700339", "588.4651556167944", "587.410328949432", "586.9501731057403", "585.8396267334707", "583.2124706812614", "584.5742851192283", "581.4196673282353")
$sigmas  = @("20.483870684032574", "23.020675159970697", "29.011885913097565", "30.051801152240728", "26.405473968522095", "22.8898502379652", "31.521944504049202", "38.983418655530954", "30.502261739669603")

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    Set-TextValue $ws.Cells.Item($row, 2) $centers[$i]
    Set-TextValue $ws.Cells.Item($row, 3) $sigmas[$i]
    Set-TextValue $ws.Cells.Item($row, 4) "lorentzian"
}

$scratch.Clear()
